# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2310
#   *_new  -> *_FV2404
# Then (re)create the worksheet table over the used range and freeze the
# header row, matching the target workbook produced by the export pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row --------------------------------------------
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J (1-10): "<Name>_old"  -> "<Name>_FV2310"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $baseNames[$i] + "_FV2310"
}

# Column K (11) stays "diff" - untouched.

# Columns L-U (12-21): "<Name>_new" -> "<Name>_FV2404"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value2 = $baseNames[$i] + "_FV2404"
}

# --- 2. Turn the used range into a real table (ListObject) --------------
$tableRange = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row --------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

"done"
